# Apply the authored changes described by the diff:
#  - Rename Sheet1 -> "Size", Sheet2 -> "Font"
#  - Move the active/selected tab from Sheet1 (Size) to Sheet2 (Font)
#  - Update each sheet's selection / active cell
#  - Reposition the workbook window (best effort)

$wb = $excel.ActiveWorkbook

# --- Rename worksheets ---------------------------------------------------
$wsSize = $wb.Worksheets.Item(1)
$wsSize.Name = "Size"

$wsFont = $wb.Worksheets.Item(2)
$wsFont.Name = "Font"

# --- Reposition / resize the workbook window (matches xWindow/yWindow in
#     the workbook.xml bookViews entry) ------------------------------------
$win = $excel.ActiveWindow
$win.Left = 5100
$win.Top = 2160
$win.Width = 21600
$win.Height = 11385

# --- Update selection on the "Size" sheet (previously tabSelected, cell H8)
$wsSize.Select()
$wsSize.Range("C11").Select()

# --- Update selection on the "Font" sheet and make it the active tab
#     (previously just a range selection A1:I11, not the active tab) ------
$wsFont.Select()
$wsFont.Range("K7").Select()

$wb.Save()
